# Evidence update for A3 sheet: record the interchain NFT transfer of
# arkNFT002 from IRISnet to Stargaze (channel 22 / channel-207 on Stargaze side).
$wb = $excel.ActiveWorkbook

$wsA2 = $wb.Worksheets.Item("A2")
$wsA3 = $wb.Worksheets.Item("A3")

# Fill in the evidence row on the "A3" sheet:
#   A2 = tx hash on IRISnet
#   B2 = ibc class-id on destination chain (Stargaze)
#   C2 = nft id
#   D2 = destination chain id (Stargaze = elgafar-1)
$wsA3.Range("A2").Value = "FB2F67ED5BD529835544C46D2C1DCB66FA6EF90954DE8A2D169FB1F67509AEC8"
$wsA3.Range("B2").Value = "wasm.stars1ve46fjrhcrum94c7d8yc2wsdz8cpuw73503e8qn9r44spr6dw0lsvmvtqh/channel-207/arkprotocol002"
$wsA3.Range("C2").Value = "arkNFT002"
$wsA3.Range("D2").Value = "elgafar-1"

# Move the active tab from "A2" to "A3", preserving each sheet's last
# selected cell.
$wsA2.Activate()
$wsA2.Range("B7").Select()

$wsA3.Activate()
$wsA3.Range("B2").Select()
